$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal (never-numeric) text value into a cell,
# using a leading apostrophe so Excel never reinterprets strings
# like "302.76" or "43.277.98" as numbers, then strip the resulting
# quote-prefix style so the cell keeps its original (default) style.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.277.98'
Set-TextValue $ws.Range("E2") '  +0.62%  '

Set-TextValue $ws.Range("D3") '2.326.00'
Set-TextValue $ws.Range("E3") '  +0.78%  '

Set-TextValue $ws.Range("E4") '  -0.10%  '

Set-TextValue $ws.Range("D5") '302.76'
Set-TextValue $ws.Range("E5") '  +0.85%  '

Set-TextValue $ws.Range("D6") '97.57'
Set-TextValue $ws.Range("E6") '  -0.29%  '

Set-TextValue $ws.Range("D7") '0.507'
Set-TextValue $ws.Range("E7") '  -1.18%  '

Set-TextValue $ws.Range("E9") '  -1.22%  '

Set-TextValue $ws.Range("D10") '35.60'
Set-TextValue $ws.Range("E10") '  -1.11%  '

Set-TextValue $ws.Range("D11") '19.53'
Set-TextValue $ws.Range("E11") '  +7.99%  '

Set-TextValue $ws.Range("D12") '0.0799'
Set-TextValue $ws.Range("E12") '  +1.10%  '

Set-TextValue $ws.Range("E13") '  +0.17%  '

Set-TextValue $ws.Range("E14") '  +2.12%  '

Set-TextValue $ws.Range("D15") '2.687.80'
Set-TextValue $ws.Range("E15") '  +0.75%  '

Set-TextValue $ws.Range("D16") '2.302.41'
Set-TextValue $ws.Range("E16") '  +0.16%  '

Set-TextValue $ws.Range("D17") '0.789'
Set-TextValue $ws.Range("E17") '  +1.00%  '

Set-TextValue $ws.Range("D18") '43.167.90'
Set-TextValue $ws.Range("E18") '  +0.51%  '

Set-TextValue $ws.Range("D19") '12.64'
Set-TextValue $ws.Range("E19") '  -1.44%  '

Set-TextValue $ws.Range("E20") '  -0.43%  '

Set-TextValue $ws.Range("D21") '6.07'
Set-TextValue $ws.Range("E21") '  +0.38%  '

Set-TextValue $ws.Range("D22") '67.89'
Set-TextValue $ws.Range("E22") '  -0.02%  '

Set-TextValue $ws.Range("D23") '237.20'
Set-TextValue $ws.Range("E23") '  +0.49%  '

Set-TextValue $ws.Range("D24") '2.25'
Set-TextValue $ws.Range("E24") '  +4.90%  '

Set-TextValue $ws.Range("E25") '  -0.05%  '

Set-TextValue $ws.Range("E26") '  +0.14%  '

Set-TextValue $ws.Range("D27") '25.01'
Set-TextValue $ws.Range("E27") '  -1.48%  '

Set-TextValue $ws.Range("E28") '  +0.99%  '

Set-TextValue $ws.Range("D29") '164.93'
Set-TextValue $ws.Range("E29") '  -0.32%  '

Set-TextValue $ws.Range("E30") '  +0.92%  '

Set-TextValue $ws.Range("D31") '33.36'
Set-TextValue $ws.Range("E31") '  -0.05%  '

Set-TextValue $ws.Range("D32") '0.999'
Set-TextValue $ws.Range("E32") '  -0.11%  '

Set-TextValue $ws.Range("D33") '17.94'
Set-TextValue $ws.Range("E33") '  +5.40%  '

Set-TextValue $ws.Range("E34") '  -0.38%  '

Set-TextValue $ws.Range("E35") '  -7.15%  '

Set-TextValue $ws.Range("E36") '  +1.38%  '

Set-TextValue $ws.Range("E37") '  -1.71%  '

Set-TextValue $ws.Range("E38") '  +0.01%  '

Set-TextValue $ws.Range("D39") '2.80'
Set-TextValue $ws.Range("E39") '  +1.94%  '

Set-TextValue $ws.Range("E40") '  +0.26%  '

Set-TextValue $ws.Range("E41") '  -0.49%  '

Set-TextValue $ws.Range("D42") '1.990.14'
Set-TextValue $ws.Range("E42") '  -1.25%  '

Set-TextValue $ws.Range("D43") '10.73'
Set-TextValue $ws.Range("E43") '  +6.50%  '

Set-TextValue $ws.Range("D44") '0.0281'
Set-TextValue $ws.Range("E44") '  -0.19%  '

Set-TextValue $ws.Range("D45") '18.22'
Set-TextValue $ws.Range("E45") '  +3.53%  '

Set-TextValue $ws.Range("D46") '2.06'
Set-TextValue $ws.Range("E46") '  -3.68%  '

Set-TextValue $ws.Range("D47") '2.79'
Set-TextValue $ws.Range("E47") '  -0.13%  '

Set-TextValue $ws.Range("B48") 'RocketPoolETH'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range("D48") '2.554.18'
Set-TextValue $ws.Range("E48") '  +0.71%  '

Set-TextValue $ws.Range("B49") 'HuobiToken'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D49") '2.87'
Set-TextValue $ws.Range("E49") '  -3.40%  '

Set-TextValue $ws.Range("D50") '53.84'
Set-TextValue $ws.Range("E50") '  -0.11%  '

Set-TextValue $ws.Range("D51") '72.16'
Set-TextValue $ws.Range("E51") '  -0.24%  '
